# "checked level 3 skill points"
# Updates several ability-rule descriptions on the Level 3 character-card
# sheet and fills in the "Requirements" (skill point) column for rows
# 11-19, which had been left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wording tweaks to existing ability text ------------------------------
$ws.Range("E8").Value  = "Target 5 by 5 foot area within 45 feet and jump to that location"
$ws.Range("E12").Value = "Destroy target summon you control. Heal for X life. Return two cards from discard."
$ws.Range("E13").Value = "Roll influence if target is unwilling. Teleport target to a 5 by 5 surface. Exhaust a card."
$ws.Range("K13").Value = "Mark 5 by 5 area. Teleport to target spot no matter the location"
$ws.Range("K14").Value = "Discard a card. Create a 10 by 10 by 10 box over a space. It has 50 hitpoints."
$ws.Range("E17").Value = "Increase either an attack die or defense die of one of your equipped weapons, and then decrease an opposite die on thatr weapon. (d4 -> d2, d2 -> 0)"
$ws.Range("F18").Value = "Exhaust"

# --- Newly checked skill-point requirements (column N) --------------------
$ws.Range("N11").Value = "6 Spirituality, 1 Finesse"
$ws.Range("N12").Value = "4 Spirituality, 2 Vigor, 1 Knowledge"
$ws.Range("N13").Value = "5 Spirituality, 2 Charisma"
$ws.Range("N14").Value = "7 Spirituality"
$ws.Range("N15").Value = "7 Charisma"
$ws.Range("N16").Value = "1 Craftsmanship, 6 Vigor"
$ws.Range("N17").Value = "7 Craftsmanship"
$ws.Range("N18").Value = "4 Perception, 1 Finesse, 2 Knowledge"
$ws.Range("N19").Value = "4 Spirituality, 2 Knowledge, 1 Vigor"

# --- Scroll position / selection left by the author after editing ---------
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("K22").Select()
